# Update crypto "Price" column (D) with refreshed quotes.
# Source cells are plain text (inline strings), so we force the
# NumberFormat to Text ("@") before writing each value — this keeps
# Excel from re-interpreting the numeric-looking string as a Number,
# matching the original text-cell representation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "249.04"
    "D3"  = "22.45"
    "D4"  = "5.343"
    "D5"  = "0.05690"
    "D6"  = "3.404"
    "D7"  = "6.312"
    "D8"  = "0.8123"
    "D9"  = "0.9202"
    "D10" = "0.1410"
    "D11" = "0.07402"
    "D13" = "0.03013"
    "D14" = "0.09360"
    "D15" = "3.712"
    "D16" = "0.001576"
    "D19" = "0.0005788"
    "D20" = "0.006434"
    "D22" = "0.001025"
    "D25" = "2.136"
    "D26" = "0.3251"
    "D41" = "0.006888"
    "D42" = "0.1064"
    "D43" = "0.002709"
    "D44" = "0.007472"
    "D48" = "0.2177"
    "D49" = "0.00002099"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
